$d = $word.ActiveDocument

# Update the date heading
[void]$d.Content.Find.Execute("2024-07-13 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-07-14 Sunday", 2)

# Update the division-fact table cells by absolute (row, column) position so that
# duplicate / cross-colliding values cannot be mismatched by a text search.
$t = $d.Tables.Item(1)

$values = @(
    @(1, 1, "19÷9=2, 1"),
    @(1, 2, "87÷5=17, 2"),
    @(1, 3, "31÷3=10, 1"),
    @(1, 4, "59÷7=8, 3"),
    @(1, 5, "67÷9=7, 4"),

    @(5, 1, "95÷7=13, 4"),
    @(5, 2, "93÷8=11, 5"),
    @(5, 3, "39÷9=4, 3"),
    @(5, 4, "47÷6=7, 5"),
    @(5, 5, "94÷4=23, 2"),

    @(9, 1, "35÷2=17, 1"),
    @(9, 2, "77÷7=11, 0"),
    @(9, 3, "16÷3=5, 1"),
    @(9, 4, "17÷8=2, 1"),
    @(9, 5, "23÷8=2, 7"),

    @(13, 1, "79÷8=9, 7"),
    @(13, 2, "65÷4=16, 1"),
    @(13, 3, "20÷5=4, 0"),
    @(13, 4, "39÷9=4, 3"),
    @(13, 5, "39÷6=6, 3"),

    @(17, 1, "84÷3=28, 0"),
    @(17, 2, "91÷8=11, 3"),
    @(17, 3, "64÷6=10, 4"),
    @(17, 4, "56÷7=8, 0"),
    @(17, 5, "49÷3=16, 1")
)

foreach ($entry in $values) {
    $row = $entry[0]
    $col = $entry[1]
    $text = $entry[2]
    [void]($t.Cell($row, $col).Range.Text = $text)
}
